$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Group 15) updated values
$ws.Range("B2").Value = -0.999490387984652
$ws.Range("C2").Value = -167965945.557761
$ws.Range("D2").Value = 0.926358919415356
$ws.Range("E2").Value = 0.99949004866267
$ws.Range("F2").Value = -0.926364242460019
$ws.Range("G2").Value = 100755.249805733
$ws.Range("H2").Value = 7053961.58094114
$ws.Range("I2").Value = -2399154.30633273
$ws.Range("J2").Value = 0.929451068167082
$ws.Range("K2").Value = 0.999621948636469
$ws.Range("L2").Value = -0.926364242460019
$ws.Range("M2").Value = 42.0068519308034
$ws.Range("N2").Value = 2931.5385345074
$ws.Range("O2").Value = -2399154.30633273

# Row 4 (Group 17) updated values
$ws.Range("B4").Value = -0.999490392340614
$ws.Range("C4").Value = -167966143.768117
$ws.Range("D4").Value = 0.92636449245941
$ws.Range("E4").Value = 0.999490393066385
$ws.Range("F4").Value = -0.926364327241407
$ws.Range("G4").Value = 100764.986438813
$ws.Range("H4").Value = 7054603.24085416
$ws.Range("I4").Value = -2399157.34660164
$ws.Range("J4").Value = 0.926397426602952
$ws.Range("K4").Value = 0.999490528013075
$ws.Range("L4").Value = -0.926364327241407
$ws.Range("M4").Value = 41.0581819232453
$ws.Range("N4").Value = 2874.40047728179
$ws.Range("O4").Value = -2399157.34660164
